$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New column widths for columns J..O (10..15)
# ---------------------------------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 9.666666666666666
$ws.Columns.Item(11).ColumnWidth = 6
$ws.Columns.Item(12).ColumnWidth = 12.833333333333334
$ws.Columns.Item(13).ColumnWidth = 14
$ws.Columns.Item(14).ColumnWidth = 9.666666666666666
$ws.Columns.Item(15).ColumnWidth = 10.166666666666666

# ---------------------------------------------------------------------
# 2) Row 1 header band: ANN (B1:E1), SVM (G1:J1), SVM-Aroon10 (L1:M1)
#    N1:O1 stay blank but pick up a "touched alignment" style (s=1)
# ---------------------------------------------------------------------
$ws.Range("N1:O1").HorizontalAlignment = 1

$ws.Range("G1").Value = "SVM"
$ws.Range("G1:J1").Merge()
$ws.Range("G1:J1").HorizontalAlignment = -4108

$ws.Range("L1").Value = "SVM - Đổi chu kỳ cho Aroon = 10"
$ws.Range("L1:M1").Merge()
$ws.Range("L1:M1").HorizontalAlignment = -4108

$ws.Range("A10").Value = "K-SVMeans"
$ws.Range("A10:J10").Merge()
$ws.Range("A10:J10").HorizontalAlignment = -4108

$ws.Range("B1").Value = "ANN"
$ws.Range("B1:E1").Merge()
$ws.Range("B1:E1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 3) New L column (second ANN-period-less block) next to the SVM table
# ---------------------------------------------------------------------
$ws.Range("L2").Value = "period = 1"
$ws.Range("L3").Value = 60.95
$ws.Range("L4").Value = 56.98
$ws.Range("L5").Value = 70.21
$ws.Range("L6").Value = 56.21
$ws.Range("L7").Value = 57.75
$ws.Range("L8").Formula = "=AVERAGE(L3:L7)"

# ---------------------------------------------------------------------
# 4) Row 8 (averages) + the new row 18 averages get a red font
# ---------------------------------------------------------------------
$ws.Range("B8,C8,G8,H8").Font.Color = 255

# ---------------------------------------------------------------------
# 5) Second table: K-SVMeans comparison (rows 11-18)
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "K = 2, Aroon 5, Volume"
$ws.Range("E11").Value = "K=2 Aroon = 5, No Volume"

$ws.Range("B12").Value = "period = 1"
$ws.Range("C12").Value = "period = 5"
$ws.Range("E12").Value = "period = 1"
$ws.Range("F12").Value = "period = 5"

$ws.Range("A13").Value = "BT6"
$ws.Range("B13").Value = 59.29
$ws.Range("C13").Value = 61.33
$ws.Range("E13").Value = 66.43

$ws.Range("A14").Value = "DHG"
$ws.Range("B14").Value = 56.98
$ws.Range("C14").Value = 56.75

$ws.Range("A15").Value = "FPT"
$ws.Range("B15").Value = 70.74
$ws.Range("C15").Value = 42.24

$ws.Range("A16").Value = "VIS"
$ws.Range("B16").Value = 46.48
$ws.Range("C16").Value = 41.84

$ws.Range("A17").Value = "VNM"
$ws.Range("B17").Value = 58.62
$ws.Range("C17").Value = 54.11

$ws.Range("A18").Value = "Total"
$ws.Range("B18").Formula = "=AVERAGE(B13:B17)"
$ws.Range("C18").Formula = "=AVERAGE(C13:C17)"
$ws.Range("B18,C18").Font.Color = 255

# Highlight a few notable cells with the theme accent color (green - theme 7)
$ws.Range("C13").Font.ThemeColor = 7
$ws.Range("B14").Font.ThemeColor = 7
$ws.Range("B17").Font.ThemeColor = 7

# ---------------------------------------------------------------------
# 6) Page setup / view cosmetics
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

$ws.Range("E18").Select()

$wb.Save()
